$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row with machine-friendly column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of Spanish connector words (de/del/el/la/los -> De/Del/El/La/Los)
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A25").Value = "Coahuila De Zaragoza"
$ws.Range("A30").Value = "Estado De México"
$ws.Range("B31").Value = "Ecatepec De Morelos"
$ws.Range("B33").Value = "San Felipe Del Progreso"
$ws.Range("B40").Value = "Acapulco De Juárez"
$ws.Range("B46").Value = "Tulancingo De Bravo"
$ws.Range("B55").Value = "La Manzanilla De La Paz"
$ws.Range("B59").Value = "Tizapán El Alto"
$ws.Range("B60").Value = "Tlajomulco De Zúñiga"
$ws.Range("B63").Value = "Zapotlán El Grande"
$ws.Range("A65").Value = "Michoacán De Ocampo"
$ws.Range("B66").Value = "Cojumatlán De Régules"
$ws.Range("B75").Value = "Ixtlán Del Río"
$ws.Range("B79").Value = "San Nicolás De Los Garza"
$ws.Range("B81").Value = "Chalcatongo De Hidalgo"
$ws.Range("B82").Value = "Oaxaca De Juárez"
$ws.Range("B88").Value = "Los Reyes De Juárez"
$ws.Range("B93").Value = "Tetela De Ocampo"
$ws.Range("B95").Value = "Cadereyta De Montes"
$ws.Range("A102").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B107").Value = "Nochistlán De Mejía"

# Fix all-caps TOTAL label to title case
$ws.Range("A110").Value = "Total"
